$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 7-12 (Feb 10/11, Feb 16/17, Mar 1 entries) entirely,
# shrinking the used range down to A1:G6.
$ws.Range("A7:G12").EntireRow.Delete()

# Rewrite rows 2-6 with the new Feb 23-27, 2026 availability/booking data.
# Column A keeps its existing date number format (style s="2"), so writing
# the 1900-date-system serial number renders as a proper date, same as
# Excel does when you type a date into an already-formatted cell.
$ws.Range("A2").Value = 46076
$ws.Range("B2").Value = "Limited"
$ws.Range("C2").Value = 18500
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Peak Season"

$ws.Range("A3").Value = 46077
$ws.Range("B3").Value = "Limited"
$ws.Range("C3").Value = 18500
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "Peak Season"

$ws.Range("A4").Value = 46078
$ws.Range("B4").Value = "Closed"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "Maintenance"

$ws.Range("A5").Value = 46079
$ws.Range("B5").Value = "Closed"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "Maintenance"

$ws.Range("A6").Value = 46080
$ws.Range("B6").Value = "Closed"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "Peak Season"

# Christmas/New Year rows are now gone from the sheet, so the corresponding
# shared strings ("Christmas Eve", "Christmas Day", "New Year's Eve",
# "New Year's Day") are no longer referenced and drop out of the saved file.

# Move the active selection to just past the new data, like a user who had
# clicked into the next empty row after editing.
$ws.Range("E7").Select()
